$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Allow WaterOilGas construction for two-phase problems ---
# The single "kroend" column is split into two columns: "krogend" (oil
# end-point relperm in oil-gas system) and "krowend" (oil end-point
# relperm in oil-water system). Insert a new column before the old
# "swirr" column (V) -- this shifts swirr/a/b/poro_ref/perm_ref/drho one
# column to the right and leaves a blank column at V for "krowend".
$ws.Columns("V").Insert()

# Rename the header that used to say "kroend" to "krogend", and label
# the newly inserted column "krowend".
$ws.Range("U1").Value = "krogend"
$ws.Range("V1").Value = "krowend"

# For every data row, the new krowend column starts out as a copy of
# the (unchanged) krogend value -- both end-points were previously a
# single "kroend" value.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 22).Value = $ws.Cells.Item($r, 21).Value()
}

# Update the selected/active cell to reflect where editing left off.
$null = $ws.Range("T7").Select()
